$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efemp1"
$ws.Cells.Item(2,3).Value = "Egfr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.536358
$ws.Cells.Item(2,8).Value = 7.609074
$ws.Cells.Item(2,9).Value = 0.04139948507354423
$ws.Cells.Item(2,10).Value = 0.04139948507354423
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.57413
$ws.Cells.Item(2,14).Value = 4.72239
$ws.Cells.Item(2,15).Value = 0.02024862668342525
$ws.Cells.Item(2,16).Value = 0.02024862668342525
$ws.Cells.Item(2,17).Value = 3.99255721854
$ws.Cells.Item(2,18).Value = 35.93301496686
$ws.Cells.Item(2,19).Value = 0.0008382827181402328
$ws.Cells.Item(2,20).Value = 0.0008382827181402328

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efemp1"
$ws.Cells.Item(3,3).Value = "Egfr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.536358
$ws.Cells.Item(3,8).Value = 7.609074
$ws.Cells.Item(3,9).Value = 0.04139948507354423
$ws.Cells.Item(3,10).Value = 0.04139948507354423
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 62.503947
$ws.Cells.Item(3,14).Value = 187.511841
$ws.Cells.Item(3,15).Value = 0.804011796385049
$ws.Cells.Item(3,16).Value = 0.8040117963850492
$ws.Cells.Item(3,17).Value = 158.532386005026
$ws.Cells.Item(3,18).Value = 1426.791474045234
$ws.Cells.Item(3,19).Value = 0.03328567436339631
$ws.Cells.Item(3,20).Value = 0.03328567436339632

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Efemp1"
$ws.Cells.Item(4,3).Value = "Egfr"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.536358
$ws.Cells.Item(4,8).Value = 7.609074
$ws.Cells.Item(4,9).Value = 0.04139948507354423
$ws.Cells.Item(4,10).Value = 0.04139948507354423
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.182903
$ws.Cells.Item(4,14).Value = 0.548709
$ws.Cells.Item(4,15).Value = 0.002352750132631058
$ws.Cells.Item(4,16).Value = 0.002352750132631058
$ws.Cells.Item(4,17).Value = 0.463907487274
$ws.Cells.Item(4,18).Value = 4.175167385466
$ws.Cells.Item(4,19).Value = 0.00009740264399763869
$ws.Cells.Item(4,20).Value = 0.00009740264399763869

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Efemp1"
$ws.Cells.Item(5,3).Value = "Egfr"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.536358
$ws.Cells.Item(5,8).Value = 7.609074
$ws.Cells.Item(5,9).Value = 0.04139948507354423
$ws.Cells.Item(5,10).Value = 0.04139948507354423
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.09380766666666666
$ws.Cells.Item(5,14).Value = 0.281423
$ws.Cells.Item(5,15).Value = 0.001206683324996365
$ws.Cells.Item(5,16).Value = 0.001206683324996365
$ws.Cells.Item(5,17).Value = 0.2379298258113333
$ws.Cells.Item(5,18).Value = 2.141368432302
$ws.Cells.Item(5,19).Value = 0.00004995606830168172
$ws.Cells.Item(5,20).Value = 0.00004995606830168172

# Row 6
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Efemp1"
$ws.Cells.Item(6,3).Value = "Egfr"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.536358
$ws.Cells.Item(6,8).Value = 7.609074
$ws.Cells.Item(6,9).Value = 0.04139948507354423
$ws.Cells.Item(6,10).Value = 0.04139948507354423
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 13.38529933333333
$ws.Cells.Item(6,14).Value = 40.155898
$ws.Cells.Item(6,15).Value = 0.1721801434738983
$ws.Cells.Item(6,16).Value = 0.1721801434738983
$ws.Cells.Item(6,17).Value = 33.94991104649467
$ws.Cells.Item(6,18).Value = 305.549199418452
$ws.Cells.Item(6,19).Value = 0.007128169279708355
$ws.Cells.Item(6,20).Value = 0.007128169279708355

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Efemp1"
$ws.Cells.Item(7,3).Value = "Egfr"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 58.41791133333334
$ws.Cells.Item(7,8).Value = 175.253734
$ws.Cells.Item(7,9).Value = 0.953521327932399
$ws.Cells.Item(7,10).Value = 0.953521327932399
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.57413
$ws.Cells.Item(7,14).Value = 4.72239
$ws.Cells.Item(7,15).Value = 0.02024862668342525
$ws.Cells.Item(7,16).Value = 0.02024862668342525
$ws.Cells.Item(7,17).Value = 91.95738676714001
$ws.Cells.Item(7,18).Value = 827.61648090426
$ws.Cells.Item(7,19).Value = 0.01930749740398705
$ws.Cells.Item(7,20).Value = 0.01930749740398705

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Efemp1"
$ws.Cells.Item(8,3).Value = "Egfr"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 58.41791133333334
$ws.Cells.Item(8,8).Value = 175.253734
$ws.Cells.Item(8,9).Value = 0.953521327932399
$ws.Cells.Item(8,10).Value = 0.953521327932399
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 62.503947
$ws.Cells.Item(8,14).Value = 187.511841
$ws.Cells.Item(8,15).Value = 0.804011796385049
$ws.Cells.Item(8,16).Value = 0.8040117963850492
$ws.Cells.Item(8,17).Value = 3651.350033829367
$ws.Cells.Item(8,18).Value = 32862.1503044643
$ws.Cells.Item(8,19).Value = 0.7666423957623856
$ws.Cells.Item(8,20).Value = 0.7666423957623857

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Efemp1"
$ws.Cells.Item(9,3).Value = "Egfr"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 58.41791133333334
$ws.Cells.Item(9,8).Value = 175.253734
$ws.Cells.Item(9,9).Value = 0.953521327932399
$ws.Cells.Item(9,10).Value = 0.953521327932399
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.182903
$ws.Cells.Item(9,14).Value = 0.548709
$ws.Cells.Item(9,15).Value = 0.002352750132631058
$ws.Cells.Item(9,16).Value = 0.002352750132631058
$ws.Cells.Item(9,17).Value = 10.68481123660067
$ws.Cells.Item(9,18).Value = 96.16330112940601
$ws.Cells.Item(9,19).Value = 0.002243397430759494
$ws.Cells.Item(9,20).Value = 0.002243397430759494

# Row 10
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Efemp1"
$ws.Cells.Item(10,3).Value = "Egfr"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 58.41791133333334
$ws.Cells.Item(10,8).Value = 175.253734
$ws.Cells.Item(10,9).Value = 0.953521327932399
$ws.Cells.Item(10,10).Value = 0.953521327932399
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.09380766666666666
$ws.Cells.Item(10,14).Value = 0.281423
$ws.Cells.Item(10,15).Value = 0.001206683324996365
$ws.Cells.Item(10,16).Value = 0.001206683324996365
$ws.Cells.Item(10,17).Value = 5.480047953720223
$ws.Cells.Item(10,18).Value = 49.320431583482
$ws.Cells.Item(10,19).Value = 0.001150598286444416
$ws.Cells.Item(10,20).Value = 0.001150598286444416

# Row 11
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Efemp1"
$ws.Cells.Item(11,3).Value = "Egfr"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 58.41791133333334
$ws.Cells.Item(11,8).Value = 175.253734
$ws.Cells.Item(11,9).Value = 0.953521327932399
$ws.Cells.Item(11,10).Value = 0.953521327932399
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 13.38529933333333
$ws.Cells.Item(11,14).Value = 40.155898
$ws.Cells.Item(11,15).Value = 0.1721801434738983
$ws.Cells.Item(11,16).Value = 0.1721801434738983
$ws.Cells.Item(11,17).Value = 781.9412296247925
$ws.Cells.Item(11,18).Value = 7037.471066623132
$ws.Cells.Item(11,19).Value = 0.1641774390488225
$ws.Cells.Item(11,20).Value = 0.1641774390488225

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Efemp1"
$ws.Cells.Item(12,3).Value = "Egfr"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.3111786666666667
$ws.Cells.Item(12,8).Value = 0.9335359999999999
$ws.Cells.Item(12,9).Value = 0.005079186994056857
$ws.Cells.Item(12,10).Value = 0.005079186994056856
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.57413
$ws.Cells.Item(12,14).Value = 4.72239
$ws.Cells.Item(12,15).Value = 0.02024862668342525
$ws.Cells.Item(12,16).Value = 0.02024862668342525
$ws.Cells.Item(12,17).Value = 0.48983567456
$ws.Cells.Item(12,18).Value = 4.408521071039999
$ws.Cells.Item(12,19).Value = 0.0001028465612979661
$ws.Cells.Item(12,20).Value = 0.0001028465612979661

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Efemp1"
$ws.Cells.Item(13,3).Value = "Egfr"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.3111786666666667
$ws.Cells.Item(13,8).Value = 0.9335359999999999
$ws.Cells.Item(13,9).Value = 0.005079186994056857
$ws.Cells.Item(13,10).Value = 0.005079186994056856
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 62.503947
$ws.Cells.Item(13,14).Value = 187.511841
$ws.Cells.Item(13,15).Value = 0.804011796385049
$ws.Cells.Item(13,16).Value = 0.8040117963850492
$ws.Cells.Item(13,17).Value = 19.449894888864
$ws.Cells.Item(13,18).Value = 175.049053999776
$ws.Cells.Item(13,19).Value = 0.004083726259267231
$ws.Cells.Item(13,20).Value = 0.004083726259267231

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Efemp1"
$ws.Cells.Item(14,3).Value = "Egfr"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.3111786666666667
$ws.Cells.Item(14,8).Value = 0.9335359999999999
$ws.Cells.Item(14,9).Value = 0.005079186994056857
$ws.Cells.Item(14,10).Value = 0.005079186994056856
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.182903
$ws.Cells.Item(14,14).Value = 0.548709
$ws.Cells.Item(14,15).Value = 0.002352750132631058
$ws.Cells.Item(14,16).Value = 0.002352750132631058
$ws.Cells.Item(14,17).Value = 0.05691551166933333
$ws.Cells.Item(14,18).Value = 0.512239605024
$ws.Cells.Item(14,19).Value = 0.00001195005787392522
$ws.Cells.Item(14,20).Value = 0.00001195005787392522

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Efemp1"
$ws.Cells.Item(15,3).Value = "Egfr"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.3111786666666667
$ws.Cells.Item(15,8).Value = 0.9335359999999999
$ws.Cells.Item(15,9).Value = 0.005079186994056857
$ws.Cells.Item(15,10).Value = 0.005079186994056856
$ws.Cells.Item(15,11).Value = 1
$ws.Cells.Item(15,12).Value = 0.3333333333333333
$ws.Cells.Item(15,13).Value = 0.09380766666666666
$ws.Cells.Item(15,14).Value = 0.281423
$ws.Cells.Item(15,15).Value = 0.001206683324996365
$ws.Cells.Item(15,16).Value = 0.001206683324996365
$ws.Cells.Item(15,17).Value = 0.02919094463644444
$ws.Cells.Item(15,18).Value = 0.2627185017279999
$ws.Cells.Item(15,19).Value = 0.000006128970250266819
$ws.Cells.Item(15,20).Value = 0.000006128970250266818

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Efemp1"
$ws.Cells.Item(16,3).Value = "Egfr"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.3111786666666667
$ws.Cells.Item(16,8).Value = 0.9335359999999999
$ws.Cells.Item(16,9).Value = 0.005079186994056857
$ws.Cells.Item(16,10).Value = 0.005079186994056856
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 13.38529933333333
$ws.Cells.Item(16,14).Value = 40.155898
$ws.Cells.Item(16,15).Value = 0.1721801434738983
$ws.Cells.Item(16,16).Value = 0.1721801434738983
$ws.Cells.Item(16,17).Value = 4.165219599480889
$ws.Cells.Item(16,18).Value = 37.486976395328
$ws.Cells.Item(16,19).Value = 0.0008745351453674677
$ws.Cells.Item(16,20).Value = 0.0008745351453674676
